# "Fruta / hortaliza, semanal" — weekly refresh of the daily price sheet.
# A new week's reading is prepended as the new row 2 (pushing the existing
# rows 2-7 down to 3-8, each keeping its own data unchanged), and the
# previously-last row (old row 7) is now duplicated into the newly
# revealed row 8 slot as part of that shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 2..7 down to 3..8, opening up a blank row 2 for the new reading.
$ws.Range("A2:T2").Insert(-4121)   # xlShiftDown

# The Insert() above copies formatting from the row above (the bold header
# row). Reset the newly opened row back to the plain/default style used by
# every other data row.
$ws.Range("A2:T2").Style = "Normal"

# Fill in the new week's reading.
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(2, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(2, 4).Value = 45043
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value = 15
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100107
$ws.Cells.Item(2, 8).Value = "Otros"
$ws.Cells.Item(2, 9).Value = 100107001
$ws.Cells.Item(2, 10).Value = "Caqui"
$ws.Cells.Item(2, 11).Value = "Fuyu"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 300
$ws.Cells.Item(2, 14).Value = 25000
$ws.Cells.Item(2, 15).Value = 26000
$ws.Cells.Item(2, 16).Value = 25500
$ws.Cells.Item(2, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(2, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(2, 19).Value = 1417
$ws.Cells.Item(2, 20).Value = 18
